# Week3 Threads, Executors, and Runnables.pptx - apply commit edits
#
# 1) Delete slide 4 ("Summary of Collection Interfaces in Java", sldId 264)
# 2) Rewrite slide 3's content placeholder text (Collection blurb -> Executors
#    blurb) and resize/move that placeholder shape.

$p = $ppt.ActivePresentation

# --- 1) Delete the "Summary of Collection Interfaces in Java" slide -------
$p.Slides.Item(4).Delete()

# --- 2) Update the content placeholder on slide 3 --------------------------
$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item(2)

# Reposition / resize the placeholder shape to its new frame.
$shp.Left = 54.0
$shp.Top = 193.2631496062992
$shp.Width = 852.0
$shp.Height = 316.86025

$tf = $shp.TextFrame
$tr = $tf.TextRange

$para1 = "Java thread pool manages the collection of Runnable threads. The worker threads execute Runnable threads from the queue. java.util.concurrent.Executors provide factory and support methods for java.util.concurrent.Executor interface to create the thread pool in java."
$para3 = "Executors is a utility class that also provides useful methods to work with ExecutorService, ScheduledExecutorService, ThreadFactory, and Callable classes through various factory methods."

# Four paragraphs: text, blank, text, blank (trailing).
$tr.Text = "$para1`r`r$para3`r"

# Bold the two java.util.concurrent.Executor(s) mentions in paragraph 1.
$p1 = $tr.Paragraphs(1, 1)
$p1.Characters(122, 30).Font.Bold = $true
$p1.Characters(193, 29).Font.Bold = $true

# Paragraph formatting: 100% line spacing + no bullet on paragraphs 1-3,
# no-bullet only on the trailing blank paragraph 4.
for ($i = 1; $i -le 3; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $para.ParagraphFormat.SpaceWithin = 1000
    $para.ParagraphFormat.Bullet.Visible = $false
}
$tr.Paragraphs(4, 1).ParagraphFormat.Bullet.Visible = $false
